$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": day_night combo strings used by the HLOOKUP table ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("D11").Value = "SaD,FaP,SaP,FaD,WaD,RaD,RaP,WaP"
$wsEv.Range("D12").Value = "FaP,SaP,FaN,SaN,WaN,RaP,WaP,RaN"

# --- Sheet "ts_12": rename ELC -> Elec in the AG column (rows 11-22) ---
$wsTs = $wb.Worksheets.Item("ts_12")
$wsTs.Range("AG11").Value = "Elec"
$wsTs.Range("AG12").Value = "Elec"
$wsTs.Range("AG13").Value = "Elec"
$wsTs.Range("AG14").Value = "Elec"
$wsTs.Range("AG15").Value = "Elec"
$wsTs.Range("AG16").Value = "Elec"
$wsTs.Range("AG17").Value = "Elec"
$wsTs.Range("AG18").Value = "Elec"
$wsTs.Range("AG19").Value = "Elec"
$wsTs.Range("AG20").Value = "Elec"
$wsTs.Range("AG21").Value = "Elec"
$wsTs.Range("AG22").Value = "Elec"

# --- Sheet "ts_12": swap the AK/AL values between rows 11 and 14 ---
$wsTs.Range("AK11").Value = "R"
$wsTs.Range("AK14").Value = "S"

$wsTs.Range("AL11").Value = 0.30301943544655252
$wsTs.Range("AL12").Value = 0.22555529847292924
$wsTs.Range("AL14").Value = 0.40439611291068944
